$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two species records that were previously on rows 7 and 8 have swapped
# places (all their field values moved to the other row). Swap the full
# rows (columns A:AY) cell by cell, using Value2 since the Value getter in
# this runtime does not reliably return the underlying data.
$lastCol = 51  # column AY

for ($col = 1; $col -le $lastCol; $col++) {
    $cell7 = $ws.Cells.Item(7, $col)
    $cell8 = $ws.Cells.Item(8, $col)

    $val7 = $cell7.Value2
    $val8 = $cell8.Value2

    $cell7.Value2 = $val8
    $cell8.Value2 = $val7
}
